$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 837.3200000000001
$ws.Range("I107").Value = 416.15384
$ws.Range("J107").Value = 1293.5834
$ws.Range("K107").Value = 416.15384
$ws.Range("L107").Value = 1293.5834
$ws.Range("M107").Value = 1503.84616
$ws.Range("N107").Value = -5133.5834

$ws.Range("H139").Value = 38182.5
$ws.Range("I139").Value = 20000
$ws.Range("J139").Value = 40780
$ws.Range("K139").Value = 20000
$ws.Range("L139").Value = 40780
$ws.Range("M139").Value = -14860
$ws.Range("N139").Value = -51060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2203.3713
$ws.Range("I61").Value = 1957.3077
$ws.Range("J61").Value = 2914.2222
$ws.Range("K61").Value = 1957.3077
$ws.Range("L61").Value = 2914.2222
$ws.Range("M61").Value = -1745.3077
$ws.Range("N61").Value = -3338.2222

$ws.Range("H64").Value = 20000
$ws.Range("J64").Value = 20000
$ws.Range("L64").Value = 20000
$ws.Range("N64").Value = -20496

$ws.Range("H67").Value = 20000
$ws.Range("J67").Value = 20000
$ws.Range("L67").Value = 20000
$ws.Range("N67").Value = -21716

$ws.Range("H122").Value = 9617102
$ws.Range("I122").Value = 20834816
$ws.Range("J122").Value = 1918.2858
$ws.Range("K122").Value = 62504448
$ws.Range("L122").Value = 5754.857400000001
$ws.Range("M122").Value = -62501998
$ws.Range("N122").Value = -10654.8574

$ws.Range("H136").Value = 2203.3713
$ws.Range("I136").Value = 1957.3077
$ws.Range("J136").Value = 2914.2222
$ws.Range("K136").Value = 5871.9231
$ws.Range("L136").Value = 8742.6666
$ws.Range("M136").Value = -3321.9231
$ws.Range("N136").Value = -13842.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 45826.668
$ws.Range("J52").Value = 45826.668
$ws.Range("L52").Value = 45826.668
$ws.Range("N52").Value = -46414.668

$ws.Range("H58").Value = 2117.0232
$ws.Range("I58").Value = 1461.3334
$ws.Range("J58").Value = 2945.2632
$ws.Range("K58").Value = 1461.3334
$ws.Range("L58").Value = 2945.2632
$ws.Range("M58").Value = -1258.3334
$ws.Range("N58").Value = -3351.2632

$ws.Range("H132").Value = 1696.9333
$ws.Range("I132").Value = 1381.52
$ws.Range("J132").Value = 2091.2
$ws.Range("K132").Value = 4144.559999999999
$ws.Range("L132").Value = 6273.599999999999
$ws.Range("M132").Value = -1614.559999999999
$ws.Range("N132").Value = -11333.6

$ws.Range("H134").Value = 15153074
$ws.Range("I134").Value = 27779254
$ws.Range("J134").Value = 1657.8
$ws.Range("K134").Value = 83337762
$ws.Range("L134").Value = 4973.4
$ws.Range("M134").Value = -83335227
$ws.Range("N134").Value = -10043.4

$ws.Range("H136").Value = 2117.0232
$ws.Range("I136").Value = 1461.3334
$ws.Range("J136").Value = 2945.2632
$ws.Range("K136").Value = 4384.0002
$ws.Range("L136").Value = 8835.7896
$ws.Range("M136").Value = -1834.0002
$ws.Range("N136").Value = -13935.7896

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7579817.5
$ws.Range("I102").Value = 18185382
$ws.Range("J102").Value = 4414.2856
$ws.Range("K102").Value = 18185382
$ws.Range("L102").Value = 4414.2856
$ws.Range("M102").Value = -18183760
$ws.Range("N102").Value = -7658.2856

$ws.Range("H107").Value = 1370.1177
$ws.Range("I107").Value = 1326.5454
$ws.Range("J107").Value = 1450
$ws.Range("K107").Value = 1326.5454
$ws.Range("L107").Value = 1450
$ws.Range("M107").Value = 593.4546
$ws.Range("N107").Value = -5290

$ws.Range("H122").Value = 2013.3334
$ws.Range("I122").Value = 2157.6
$ws.Range("J122").Value = 1652.6666
$ws.Range("K122").Value = 6472.799999999999
$ws.Range("L122").Value = 4957.9998
$ws.Range("M122").Value = -4022.799999999999
$ws.Range("N122").Value = -9857.9998

$ws.Range("H132").Value = 13896317
$ws.Range("I132").Value = 41681412
$ws.Range("J132").Value = 3768.9375
$ws.Range("K132").Value = 125044236
$ws.Range("L132").Value = 11306.8125
$ws.Range("M132").Value = -125041706
$ws.Range("N132").Value = -16366.8125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2592.8667
$ws.Range("I7").Value = 2590.6667
$ws.Range("J7").Value = 2601.6667
$ws.Range("K7").Value = 2590.6667
$ws.Range("L7").Value = 2601.6667
$ws.Range("M7").Value = -2478.6667
$ws.Range("N7").Value = -2825.6667

$ws.Range("H82").Value = 3411.111
$ws.Range("I82").Value = 2500
$ws.Range("J82").Value = 3671.4285
$ws.Range("K82").Value = 2500
$ws.Range("L82").Value = 3671.4285
$ws.Range("M82").Value = -2139
$ws.Range("N82").Value = -4393.4285

$ws.Range("H85").Value = 3411.111
$ws.Range("I85").Value = 2500
$ws.Range("J85").Value = 3671.4285
$ws.Range("K85").Value = 2500
$ws.Range("L85").Value = 3671.4285
$ws.Range("M85").Value = -1252
$ws.Range("N85").Value = -6167.4285

$ws.Range("H122").Value = 4012.2354
$ws.Range("I122").Value = 4086.2856
$ws.Range("J122").Value = 3666.6667
$ws.Range("K122").Value = 12258.8568
$ws.Range("L122").Value = 11000.0001
$ws.Range("M122").Value = -9808.856800000001
$ws.Range("N122").Value = -15900.0001

$ws.Range("H126").Value = 2592.8667
$ws.Range("I126").Value = 2590.6667
$ws.Range("J126").Value = 2601.6667
$ws.Range("K126").Value = 7772.000100000001
$ws.Range("L126").Value = 7805.000100000001
$ws.Range("M126").Value = -5302.000100000001
$ws.Range("N126").Value = -12745.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4216.6665
$ws.Range("J62").Value = 4260
$ws.Range("L62").Value = 4260
$ws.Range("N62").Value = -5508

$ws.Range("H65").Value = 4216.6665
$ws.Range("J65").Value = 4260
$ws.Range("L65").Value = 21300
$ws.Range("N65").Value = -27540

$ws.Range("H107").Value = 125000376
$ws.Range("I107").Value = 500000000
$ws.Range("K107").Value = 1500000000
$ws.Range("M107").Value = -1499998080

$ws.Range("H126").Value = 1564.1818
$ws.Range("I126").Value = 880.375
$ws.Range("J126").Value = 1954.9286
$ws.Range("K126").Value = 2641.125
$ws.Range("L126").Value = 5864.7858
$ws.Range("M126").Value = -171.125
$ws.Range("N126").Value = -10804.7858

$ws.Range("H136").Value = 7264.7393
$ws.Range("I136").Value = 2728
$ws.Range("J136").Value = 11423.417
$ws.Range("K136").Value = 8184
$ws.Range("L136").Value = 34270.251
$ws.Range("M136").Value = -5634
$ws.Range("N136").Value = -39370.251
